$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "courses"
$ws.Range("E2").Value = "java,c"
$ws.Range("E3").Value = "js,php"
$ws.Range("E4").Value = "c++,java"

$ws.Range("E2").Select()
